$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.3268493333333333
$ws.Range("N2").Value = 0.980548
$ws.Range("O2").Value = 0.05842197836270246
$ws.Range("P2").Value = 0.05842197836270247
$ws.Range("Q2").Value = 27.12446820972933
$ws.Range("R2").Value = 244.120213887564
$ws.Range("S2").Value = 0.02622857122907955
$ws.Range("T2").Value = 0.02622857122907955

# Row 3
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.9415780216372975
$ws.Range("P3").Value = 0.9415780216372976
$ws.Range("Q3").Value = 437.1608738807403
$ws.Range("R3").Value = 3934.447864926663
$ws.Range("S3").Value = 0.4227218403137157
$ws.Range("T3").Value = 0.4227218403137157

# Row 4
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("M4").Value = 0.3268493333333333
$ws.Range("N4").Value = 0.980548
$ws.Range("O4").Value = 0.05842197836270246
$ws.Range("P4").Value = 0.05842197836270247
$ws.Range("Q4").Value = 20.63745952987378
$ws.Range("R4").Value = 185.737135768864
$ws.Range("S4").Value = 0.01995582265728565
$ws.Range("T4").Value = 0.01995582265728565

# Row 5
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("O5").Value = 0.9415780216372975
$ws.Range("P5").Value = 0.9415780216372976
$ws.Range("S5").Value = 0.3216249182993707
$ws.Range("T5").Value = 0.3216249182993707

# Row 6
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("M6").Value = 0.3268493333333333
$ws.Range("N6").Value = 0.980548
$ws.Range("O6").Value = 0.05842197836270246
$ws.Range("P6").Value = 0.05842197836270247
$ws.Range("Q6").Value = 12.65558722940533
$ws.Range("R6").Value = 113.900285064648
$ws.Range("S6").Value = 0.01223758447633727
$ws.Range("T6").Value = 0.01223758447633727

# Row 7
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("O7").Value = 0.9415780216372975
$ws.Range("P7").Value = 0.9415780216372976
$ws.Range("Q7").Value = 203.9681489754074
$ws.Range("S7").Value = 0.1972312630242112
$ws.Range("T7").Value = 0.1972312630242112
